$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text: Volume number 2 -> 3, and week date range ---
$ws.Range("A8").Value = "Volume 30   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/16/2023  Through  1/22/2023"

# --- Update crime-stat grid (rows 14-29) ---
$ws.Range("L14").Copy()
$ws.Range("M14").PasteSpecial(-4122)
$ws.Range("M14").Value = -100

$ws.Range("C14").Copy()
$ws.Range("F15").PasteSpecial(-4104)

$ws.Range("H15").Value = -100

$ws.Range("C16").Value = 3

$ws.Range("D16").Value = 2

$ws.Range("E16").Value = 50

$ws.Range("F16").Value = 9

$ws.Range("G16").Value = 13

$ws.Range("H16").Value = -30.76923076923

$ws.Range("I16").Value = 7

$ws.Range("J16").Value = 12

$ws.Range("K16").Value = -41.666666666666

$ws.Range("L16").Value = 75

$ws.Range("M16").Value = -30

$ws.Range("N16").Value = -90.90909090909

$ws.Range("C14").Copy()
$ws.Range("C17").PasteSpecial(-4104)

$ws.Range("E17").Value = -100

$ws.Range("F17").Value = 9

$ws.Range("H17").Value = 50

$ws.Range("J17").Value = 5

$ws.Range("K17").Value = 0

$ws.Range("L17").Value = -28.571428571428

$ws.Range("M17").Value = -50

$ws.Range("N17").Value = -82.142857142857

$ws.Range("C18").Value = 1

$ws.Range("D18").Value = 5

$ws.Range("E18").Value = -80

$ws.Range("F18").Value = 12

$ws.Range("H18").Value = -7.692307692307

$ws.Range("I18").Value = 8

$ws.Range("J18").Value = 11

$ws.Range("K18").Value = -27.272727272727

$ws.Range("L18").Value = -38.461538461538

$ws.Range("M18").Value = -38.461538461538

$ws.Range("N18").Value = -94.871794871794

$ws.Range("C19").Value = 32

$ws.Range("D19").Value = 23

$ws.Range("E19").Value = 39.130434782608

$ws.Range("F19").Value = 185

$ws.Range("G19").Value = 129

$ws.Range("H19").Value = 43.410852713178

$ws.Range("I19").Value = 129

$ws.Range("J19").Value = 82

$ws.Range("K19").Value = 57.317073170731

$ws.Range("L19").Value = 130.357142857143

$ws.Range("M19").Value = 27.722772277227

$ws.Range("N19").Value = -72.435897435897

$ws.Range("C20").Value = 1

$ws.Range("C14").Copy()
$ws.Range("D20").PasteSpecial(-4104)

$ws.Range("E14").Copy()
$ws.Range("E20").PasteSpecial(-4104)

$ws.Range("G20").Value = 10

$ws.Range("H20").Value = -40

$ws.Range("I20").Value = 5

$ws.Range("K20").Value = -37.5

$ws.Range("L20").Value = 150

$ws.Range("M20").Value = 400

$ws.Range("N20").Value = -84.375

$ws.Range("C21").Value = 37

$ws.Range("D21").Value = 32

$ws.Range("E21").Value = 15.625

$ws.Range("F21").Value = 221

$ws.Range("G21").Value = 172

$ws.Range("H21").Value = 28.488372093023

$ws.Range("I21").Value = 154

$ws.Range("J21").Value = 119

$ws.Range("K21").Value = 29.411764705882

$ws.Range("L21").Value = 85.542168674698

$ws.Range("M21").Value = 10

$ws.Range("N21").Value = -79.973992197659

$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4104)

$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4104)

$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4104)

$ws.Range("F22").Value = 5

$ws.Range("G22").Value = 2

$ws.Range("H22").Value = 150

$ws.Range("M22").Value = 0

$ws.Range("C24").Value = 42

$ws.Range("E24").Value = 31.25

$ws.Range("F24").Value = 201

$ws.Range("G24").Value = 152

$ws.Range("H24").Value = 32.236842105263

$ws.Range("I24").Value = 147

$ws.Range("J24").Value = 91

$ws.Range("K24").Value = 61.538461538461

$ws.Range("L24").Value = 51.546391752577

$ws.Range("M24").Value = 56.382978723404

$ws.Range("C25").Value = 8

$ws.Range("D25").Value = 3

$ws.Range("E25").Value = 166.666666666667

$ws.Range("F25").Value = 44

$ws.Range("G25").Value = 35

$ws.Range("H25").Value = 25.714285714285

$ws.Range("I25").Value = 34

$ws.Range("J25").Value = 29

$ws.Range("K25").Value = 17.241379310344

$ws.Range("L25").Value = 161.538461538462

$ws.Range("M25").Value = 3.030303030303

$ws.Range("C14").Copy()
$ws.Range("D26").PasteSpecial(-4104)

$ws.Range("E14").Copy()
$ws.Range("E26").PasteSpecial(-4104)

$ws.Range("C14").Copy()
$ws.Range("F26").PasteSpecial(-4104)

$ws.Range("H26").Value = -100

$ws.Range("G15").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1

$ws.Range("D27").Value = 2

$ws.Range("E27").Value = -50

$ws.Range("G27").Value = 8

$ws.Range("H27").Value = -50

$ws.Range("G15").Copy()
$ws.Range("I27").PasteSpecial(-4122)
$ws.Range("I27").Value = 1

$ws.Range("J27").Value = 7

$ws.Range("K27").Value = -85.714285714285

$ws.Range("L14").Copy()
$ws.Range("L27").PasteSpecial(-4122)
$ws.Range("L27").Value = -66.666666666666

$ws.Range("G15").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 1

$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4104)

$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4104)

$ws.Range("G15").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("F28").Value = 1

$ws.Range("H28").Value = 0

$ws.Range("G15").Copy()
$ws.Range("I28").PasteSpecial(-4122)
$ws.Range("I28").Value = 1

$ws.Range("K28").Value = 0

$ws.Range("N28").Value = 0

$ws.Range("G15").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = 1

$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4104)

$ws.Range("E14").Copy()
$ws.Range("E29").PasteSpecial(-4104)

$ws.Range("G15").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("F29").Value = 1

$ws.Range("H29").Value = 0

$ws.Range("G15").Copy()
$ws.Range("I29").PasteSpecial(-4122)
$ws.Range("I29").Value = 1

$ws.Range("K29").Value = 0

$ws.Range("N29").Value = 0

$excel.CutCopyMode = 0